$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (1-based, referring to the original/current row numbering at the time
# each statement executes) that must be removed entirely. We process them
# from the bottom of the sheet upward so that earlier deletions never shift
# the row number of a not-yet-deleted row out from under us.
#   2  -> 004503381 FREDERICO 1138598.99
#   3  -> 004503374 PAULA     710642.69
#   4  -> 005642649 VR        500000
#   5  -> 004525587 MARIANA   214382.25   (re-added below with the new balance)
#   6  -> 004212438 KENIA     137437.5
#   10 -> 005701765 F         30691.62
#   11 -> 004381180 HFR       24185.61
#   12 -> 004567880 LUANA     23089.5
#   14 -> 005599726 JORGE     22119.17
#   15 -> 004238436 DIEGO     12072.31
#   16 -> 005581299 ZILDA     4400.85
#   19 -> 004382902 LEILA     3292.45
#   21 -> 005338054 ELAINE    1058.99
$rowsToDelete = @(21, 19, 16, 15, 14, 12, 11, 10, 6, 5, 4, 3, 2)
foreach ($rowNum in $rowsToDelete) {
    $ws.Rows.Item($rowNum).Delete()
}

# After the deletions above, the row that used to hold CAROLINA (originally
# row 7) has shifted up to row 2. Re-insert MARIANA's row right after it,
# with her balance updated to 64382.25.
$ws.Rows.Item(3).Insert()
$ws.Cells.Item(3, 1).Value = "'004525587"
$ws.Cells.Item(3, 2).Value = "MARIANA"
$ws.Cells.Item(3, 3).Value = 64382.25
